# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect freshly scraped data (commit: Update gh-pages to output
# generated at 456a3b4).

$wb = $excel.ActiveWorkbook

$updates = @{
    "F2"  = 72
    "F4"  = 2069
    "F5"  = 356
    "F6"  = 609
    "F7"  = 98
    "F9"  = 10625
    "F12" = 282
    "F13" = 202
    "F15" = 7504
    "F17" = 717
    "F18" = 246
    "F20" = 3326
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($cell in $updates.Keys) {
    $ws1.Range($cell).Value = $updates[$cell]
}

$updates2 = @{
    "F2"  = 72
    "F4"  = 2069
    "F5"  = 356
    "F6"  = 609
    "F8"  = 98
    "F12" = 10625
    "F15" = 282
    "F16" = 202
    "F18" = 7504
    "F20" = 717
    "F21" = 246
    "F23" = 3326
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($cell in $updates2.Keys) {
    $ws4.Range($cell).Value = $updates2[$cell]
}
